$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the three-reel combo rows ("c _ _", "_cc", "c_c", "cc_") - rows
#    17-20 - which are being dropped from the payout table. Deleting shifts
#    every row below up by four and keeps all formulas/refs self-consistent
#    (shared formula range, SUM range, etc. all auto-adjust).
# ---------------------------------------------------------------------------
$ws.Rows("17:20").Delete()

# ---------------------------------------------------------------------------
# 2. Turn the old "_ _ c" / "_ c _" single-cherry rows (now rows 15 & 16)
#    into the new 1-credit / 2-credit paylines, with their own odds.
# ---------------------------------------------------------------------------
$ws.Range("A15").Value2 = "1c"
$ws.Range("B15").Value2 = 1
$ws.Range("C15").Value2 = 4.1152263374485596

$ws.Range("A16").Value2 = "2c"
$ws.Range("B16").Value2 = 2
$ws.Range("C16").Value2 = 12.345679012345679
$ws.Range("J16").Value2 = "-"

# ---------------------------------------------------------------------------
# 3. Tweak a handful of payout values further down the (now shorter) table.
# ---------------------------------------------------------------------------
$ws.Range("B22").Value2 = 35
$ws.Range("B24").Value2 = 35
$ws.Range("B25").Value2 = 35

# ---------------------------------------------------------------------------
# 4. Clear the stray "D30/10" formula that Excel's row-delete left behind in
#    column E of the (now) Total row, and the old Chance column format.
# ---------------------------------------------------------------------------
$ws.Range("E26").ClearContents()

# ---------------------------------------------------------------------------
# 5. Build the new bold, boxed "Total:" / "Return:" summary block.
#    Seed the formatting from B1 (bold, General, no border) so every cell
#    starts from a clean, already-existing style before we layer the
#    per-cell number formats and border edges on top.
# ---------------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("A26:D27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A26").Value2 = "Total:"
$ws.Range("D26").Formula = "=SUM(D14:D25)"
$ws.Range("A27").Value2 = "Return:"
$ws.Range("D27").Formula = "=D26/1000"

$ws.Range("C26").NumberFormat = "@"
$ws.Range("A27").NumberFormat = "@"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "0.00%"

# Top edge of the box (row 26): left+top corner, top-only, top+right corner.
$ws.Range("A26").Borders(7).LineStyle = 1
$ws.Range("A26").Borders(8).LineStyle = 1
$ws.Range("B26:C26").Borders(8).LineStyle = 1
$ws.Range("D26").Borders(8).LineStyle = 1
$ws.Range("D26").Borders(10).LineStyle = 1

# Bottom edge of the box (row 27): left+bottom corner, bottom-only, bottom+right corner.
$ws.Range("A27").Borders(7).LineStyle = 1
$ws.Range("A27").Borders(9).LineStyle = 1
$ws.Range("B27:C27").Borders(9).LineStyle = 1
$ws.Range("D27").Borders(9).LineStyle = 1
$ws.Range("D27").Borders(10).LineStyle = 1

# ---------------------------------------------------------------------------
# 6. Cosmetic: scroll position / selection, matching the saved view state.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("J42").Select()
